$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.556.01"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.26"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.49"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.52"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.643.97"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.85"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("E15").Value = "  +4.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.126.19"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.444.60"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.643.24"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  +3.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "368.70"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.85"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.29"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("E28").Value = "  +6.48%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.97"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("E32").Value = "  +4.23%  "
$ws.Range("E33").Value = "  +4.69%  "
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +3.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.50"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("E39").Value = "  +4.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.23"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.38"
$ws.Range("E41").Value = "  +3.47%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.367"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  +4.50%  "
$ws.Range("E44").Value = "  +3.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0321"
$ws.Range("E45").Value = "  +12.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.56"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("E49").Value = "  +0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.71"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.95"
$ws.Range("E51").Value = "  +2.10%  "
